$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B78").Value = 44975
$ws.Range("D78").Value = 282.73
$ws.Range("E78").Value = 320.6
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0

$ws.Range("B79").Value = 63798
$ws.Range("D79").Value = 282.73
$ws.Range("E79").Value = 300.58
$ws.Range("F79").Value = 1
$ws.Range("G79").Value = 282.73

$ws.Range("B82").Value = 44977
$ws.Range("D82").Value = 377.19
$ws.Range("E82").Value = 427.72
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 0

$ws.Range("B83").Value = 63799
$ws.Range("D83").Value = 377.19
$ws.Range("E83").Value = 401
$ws.Range("F83").Value = 1
$ws.Range("G83").Value = 377.19

$ws.Range("B84").Value = 63792
$ws.Range("D84").Value = 916.9400000000001
$ws.Range("E84").Value = 974.8200000000001
$ws.Range("F84").Value = 4
$ws.Range("G84").Value = 3667.76

$ws.Range("B85").Value = 44959
$ws.Range("D85").Value = 916.9400000000001
$ws.Range("E85").Value = 1039.82
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0

$ws.Range("B96").Value = 59710
$ws.Range("D96").Value = 172.04
$ws.Range("E96").Value = 205.53
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0

$ws.Range("B97").Value = 64477
$ws.Range("D97").Value = 172.04
$ws.Range("E97").Value = 182.89
$ws.Range("F97").Value = 28
$ws.Range("G97").Value = 4817.12

$ws.Range("B98").Value = 64479
$ws.Range("D98").Value = 82.14
$ws.Range("E98").Value = 87.33
$ws.Range("F98").Value = 35
$ws.Range("G98").Value = 2874.9

$ws.Range("B99").Value = 59712
$ws.Range("D99").Value = 82.14
$ws.Range("E99").Value = 98.13
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0

$ws.Range("B131").Value = 64466
$ws.Range("D131").Value = 62.17
$ws.Range("E131").Value = 66.09
$ws.Range("F131").Value = 115
$ws.Range("G131").Value = 7149.55

$ws.Range("B132").Value = 59659
$ws.Range("D132").Value = 62.17
$ws.Range("E132").Value = 74.29000000000001
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0

$ws.Range("B145").Value = 59602
$ws.Range("D145").Value = 307.77
$ws.Range("E145").Value = 349.02
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0

$ws.Range("B146").Value = 64458
$ws.Range("D146").Value = 307.77
$ws.Range("E146").Value = 327.2
$ws.Range("F146").Value = 23
$ws.Range("G146").Value = 7078.71

$ws.Range("B147").Value = 64459
$ws.Range("D147").Value = 281.55
$ws.Range("E147").Value = 299.33
$ws.Range("F147").Value = 3
$ws.Range("G147").Value = 844.65

$ws.Range("B148").Value = 59603
$ws.Range("D148").Value = 281.55
$ws.Range("E148").Value = 319.27
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0

$ws.Range("B150").Value = 64461
$ws.Range("D150").Value = 598.4299999999999
$ws.Range("E150").Value = 636.21
$ws.Range("F150").Value = 5
$ws.Range("G150").Value = 2992.15

$ws.Range("B151").Value = 59634
$ws.Range("D151").Value = 598.4299999999999
$ws.Range("E151").Value = 678.63
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0

$ws.Range("B152").Value = 64452
$ws.Range("D152").Value = 47.3
$ws.Range("E152").Value = 50.27
$ws.Range("F152").Value = 39
$ws.Range("G152").Value = 1844.7

$ws.Range("B153").Value = 59595
$ws.Range("D153").Value = 47.3
$ws.Range("E153").Value = 56.51
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0

$ws.Range("B157").Value = 64196
$ws.Range("D157").Value = 32143.58
$ws.Range("E157").Value = 38403.53
$ws.Range("F157").Value = 1
$ws.Range("G157").Value = 32143.58

$ws.Range("B158").Value = 54863
$ws.Range("D158").Value = 32143.58
$ws.Range("E158").Value = 41658.07
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0

$ws.Range("B175").Value = 53144
$ws.Range("D175").Value = 1763.16
$ws.Range("E175").Value = 1999.42
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0

$ws.Range("B176").Value = 64076
$ws.Range("D176").Value = 1763.16
$ws.Range("E176").Value = 1874.46
$ws.Range("F176").Value = 1
$ws.Range("G176").Value = 1763.16

$ws.Range("B182").Value = 63911
$ws.Range("D182").Value = 49.48
$ws.Range("E182").Value = 52.62
$ws.Range("F182").Value = 47
$ws.Range("G182").Value = 2325.56

$ws.Range("B183").Value = 48678
$ws.Range("D183").Value = 49.48
$ws.Range("E183").Value = 59.12
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0

$ws.Range("B189").Value = 63902
$ws.Range("D189").Value = 32.02
$ws.Range("E189").Value = 34.04
$ws.Range("F189").Value = 2
$ws.Range("G189").Value = 64.04000000000001

$ws.Range("B190").Value = 48654
$ws.Range("D190").Value = 32.02
$ws.Range("E190").Value = 38.26
$ws.Range("F190").Value = 1
$ws.Range("G190").Value = 32.02

$ws.Range("B225").Value = 57540
$ws.Range("D225").Value = 39.49
$ws.Range("E225").Value = 44.78
$ws.Range("F225").Value = 0
$ws.Range("G225").Value = 0

$ws.Range("B226").Value = 64325
$ws.Range("D226").Value = 39.49
$ws.Range("E226").Value = 41.98
$ws.Range("F226").Value = 24
$ws.Range("G226").Value = 947.76

$ws.Range("B251").Value = 46905
$ws.Range("D251").Value = 42.64
$ws.Range("E251").Value = 48.35
$ws.Range("F251").Value = 0
$ws.Range("G251").Value = 0

$ws.Range("B252").Value = 63848
$ws.Range("D252").Value = 42.64
$ws.Range("E252").Value = 45.33
$ws.Range("F252").Value = 17
$ws.Range("G252").Value = 724.88

$ws.Range("B253").Value = 56931
$ws.Range("D253").Value = 53.53
$ws.Range("E253").Value = 60.7
$ws.Range("F253").Value = 0
$ws.Range("G253").Value = 0

$ws.Range("B254").Value = 64285
$ws.Range("D254").Value = 53.53
$ws.Range("E254").Value = 56.92
$ws.Range("F254").Value = 58
$ws.Range("G254").Value = 3104.74

$ws.Range("B366").Value = 62997
$ws.Range("D366").Value = 305.84
$ws.Range("E366").Value = 325.16
$ws.Range("F366").Value = 72
$ws.Range("G366").Value = 22020.48

$ws.Range("B367").Value = 57854
$ws.Range("D367").Value = 305.84
$ws.Range("E367").Value = 325.16
$ws.Range("F367").Value = 2
$ws.Range("G367").Value = 611.6799999999999

$ws.Range("B370").Value = 57077
$ws.Range("D370").Value = 93.08
$ws.Range("E370").Value = 111.2
$ws.Range("F370").Value = 1
$ws.Range("G370").Value = 93.08

$ws.Range("B371").Value = 63565
$ws.Range("D371").Value = 102.71
$ws.Range("E371").Value = 109.19
$ws.Range("F371").Value = 60
$ws.Range("G371").Value = 6162.6

$ws.Range("B374").Value = 63548
$ws.Range("D374").Value = 107.09
$ws.Range("E374").Value = 113.85
$ws.Range("F374").Value = 14
$ws.Range("G374").Value = 1499.26

$ws.Range("B375").Value = 57842
$ws.Range("D375").Value = 107.09
$ws.Range("E375").Value = 127.95
$ws.Range("F375").Value = 0
$ws.Range("G375").Value = 0

$ws.Range("B404").Value = 63520
$ws.Range("D404").Value = 144.28
$ws.Range("E404").Value = 153.4
$ws.Range("F404").Value = 97
$ws.Range("G404").Value = 13995.16

$ws.Range("B405").Value = 55373
$ws.Range("D405").Value = 144.28
$ws.Range("E405").Value = 163.62
$ws.Range("F405").Value = 3
$ws.Range("G405").Value = 432.84

$ws.Range("B408").Value = 57802
$ws.Range("D408").Value = 143.48
$ws.Range("E408").Value = 162.71
$ws.Range("F408").Value = 1
$ws.Range("G408").Value = 143.48

$ws.Range("B409").Value = 63531
$ws.Range("D409").Value = 143.48
$ws.Range("E409").Value = 152.53
$ws.Range("F409").Value = 80
$ws.Range("G409").Value = 11478.4

$ws.Range("B410").Value = 63571
$ws.Range("D410").Value = 143.48
$ws.Range("E410").Value = 152.53
$ws.Range("F410").Value = 29
$ws.Range("G410").Value = 4160.92

$ws.Range("B427").Value = 63521
$ws.Range("D427").Value = 115.05
$ws.Range("E427").Value = 137.45
$ws.Range("F427").Value = 16
$ws.Range("G427").Value = 1840.8

$ws.Range("B428").Value = 55378
$ws.Range("D428").Value = 115.05
$ws.Range("E428").Value = 116.49
$ws.Range("F428").Value = 0
$ws.Range("G428").Value = 0

$ws.Range("B435").Value = 63652
$ws.Range("D435").Value = 52.13
$ws.Range("E435").Value = 55.42
$ws.Range("F435").Value = 250
$ws.Range("G435").Value = 13032.5

$ws.Range("B436").Value = 57885
$ws.Range("D436").Value = 52.13
$ws.Range("E436").Value = 62.28
$ws.Range("F436").Value = 4
$ws.Range("G436").Value = 208.52

$ws.Range("B438").Value = 61605
$ws.Range("D438").Value = 111.96
$ws.Range("E438").Value = 133.78
$ws.Range("F438").Value = 2
$ws.Range("G438").Value = 223.92

$ws.Range("B439").Value = 63563
$ws.Range("D439").Value = 111.96
$ws.Range("E439").Value = 119.04
$ws.Range("F439").Value = 15
$ws.Range("G439").Value = 1679.4

$ws.Range("B442").Value = 63564
$ws.Range("D442").Value = 129.01
$ws.Range("E442").Value = 137.16
$ws.Range("F442").Value = 57
$ws.Range("G442").Value = 7353.57

$ws.Range("B443").Value = 61608
$ws.Range("D443").Value = 129.01
$ws.Range("E443").Value = 154.12
$ws.Range("F443").Value = 1
$ws.Range("G443").Value = 129.01

$ws.Range("B445").Value = 60325
$ws.Range("D445").Value = 126.86
$ws.Range("E445").Value = 151.57
$ws.Range("F445").Value = 2
$ws.Range("G445").Value = 253.72

$ws.Range("B446").Value = 63560
$ws.Range("D446").Value = 126.86
$ws.Range("E446").Value = 134.87
$ws.Range("F446").Value = 104
$ws.Range("G446").Value = 13193.44

$ws.Range("B454").Value = 57817
$ws.Range("D454").Value = 79.81
$ws.Range("E454").Value = 95.34999999999999
$ws.Range("F454").Value = 3
$ws.Range("G454").Value = 239.43

$ws.Range("B455").Value = 62865
$ws.Range("D455").Value = 79.81
$ws.Range("E455").Value = 95.34999999999999
$ws.Range("F455").Value = 151
$ws.Range("G455").Value = 12051.31

$ws.Range("B485").Value = 57856
$ws.Range("D485").Value = 171.33
$ws.Range("E485").Value = 204.69
$ws.Range("F485").Value = 2
$ws.Range("G485").Value = 342.66

$ws.Range("B486").Value = 63007
$ws.Range("D486").Value = 171.33
$ws.Range("E486").Value = 204.69
$ws.Range("F486").Value = 984
$ws.Range("G486").Value = 168588.72

$ws.Range("B502").Value = 49690
$ws.Range("D502").Value = 153.41
$ws.Range("E502").Value = 183.29
$ws.Range("F502").Value = 0
$ws.Range("G502").Value = 0

$ws.Range("B503").Value = 63940
$ws.Range("D503").Value = 153.41
$ws.Range("E503").Value = 163.11
$ws.Range("F503").Value = 24
$ws.Range("G503").Value = 3681.84

$ws.Range("B506").Value = 57903
$ws.Range("D506").Value = 446.65
$ws.Range("E506").Value = 533.63
$ws.Range("F506").Value = 0
$ws.Range("G506").Value = 0

$ws.Range("B507").Value = 64361
$ws.Range("D507").Value = 446.65
$ws.Range("E507").Value = 474.85
$ws.Range("F507").Value = 4
$ws.Range("G507").Value = 1786.6

$ws.Range("B541").Value = 54533
$ws.Range("D541").Value = 321.01
$ws.Range("E541").Value = 364.02
$ws.Range("F541").Value = 0
$ws.Range("G541").Value = 0

$ws.Range("B542").Value = 64191
$ws.Range("D542").Value = 321.01
$ws.Range("E542").Value = 341.28
$ws.Range("F542").Value = 2
$ws.Range("G542").Value = 642.02

$ws.Range("B620").Value = 51399
$ws.Range("D620").Value = 755
$ws.Range("E620").Value = 856.1799999999999
$ws.Range("F620").Value = 0
$ws.Range("G620").Value = 0

$ws.Range("B621").Value = 64030
$ws.Range("D621").Value = 755
$ws.Range("E621").Value = 802.6799999999999
$ws.Range("F621").Value = 1
$ws.Range("G621").Value = 755

$ws.Range("B632").Value = 58047
$ws.Range("D632").Value = 105.54
$ws.Range("E632").Value = 126.1
$ws.Range("F632").Value = 54
$ws.Range("G632").Value = 5699.16

$ws.Range("B633").Value = 47097
$ws.Range("D633").Value = 112.28
$ws.Range("E633").Value = 134.16
$ws.Range("F633").Value = 15
$ws.Range("G633").Value = 1684.2

$ws.Range("B714").Value = 50911
$ws.Range("D714").Value = 159.32
$ws.Range("E714").Value = 180.67
$ws.Range("F714").Value = 0
$ws.Range("G714").Value = 0

$ws.Range("B715").Value = 64013
$ws.Range("D715").Value = 159.32
$ws.Range("E715").Value = 169.37
$ws.Range("F715").Value = 54
$ws.Range("G715").Value = 8603.280000000001

$ws.Range("B718").Value = 50910
$ws.Range("D718").Value = 159.32
$ws.Range("E718").Value = 180.67
$ws.Range("F718").Value = 0
$ws.Range("G718").Value = 0

$ws.Range("B719").Value = 64012
$ws.Range("D719").Value = 159.32
$ws.Range("E719").Value = 169.37
$ws.Range("F719").Value = 78
$ws.Range("G719").Value = 12426.96

$ws.Range("B748").Value = 64244
$ws.Range("D748").Value = 670.65
$ws.Range("E748").Value = 712.99
$ws.Range("F748").Value = 2
$ws.Range("G748").Value = 1341.3

$ws.Range("B749").Value = 55658
$ws.Range("D749").Value = 670.65
$ws.Range("E749").Value = 801.25
$ws.Range("F749").Value = 0
$ws.Range("G749").Value = 0

$ws.Range("B750").Value = 55635
$ws.Range("D750").Value = 583.95
$ws.Range("E750").Value = 697.6900000000001
$ws.Range("F750").Value = 0
$ws.Range("G750").Value = 0

$ws.Range("B751").Value = 64233
$ws.Range("D751").Value = 583.95
$ws.Range("E751").Value = 620.8099999999999
$ws.Range("F751").Value = 4
$ws.Range("G751").Value = 2335.8

$ws.Range("B752").Value = 55655
$ws.Range("D752").Value = 583.95
$ws.Range("E752").Value = 697.6900000000001
$ws.Range("F752").Value = 0
$ws.Range("G752").Value = 0

$ws.Range("B753").Value = 64243
$ws.Range("D753").Value = 583.95
$ws.Range("E753").Value = 620.8099999999999
$ws.Range("F753").Value = 28
$ws.Range("G753").Value = 16350.6

$ws.Range("B776").Value = 46270
$ws.Range("D776").Value = 6.85
$ws.Range("E776").Value = 8.199999999999999
$ws.Range("F776").Value = 0
$ws.Range("G776").Value = 0

$ws.Range("B777").Value = 63810
$ws.Range("D777").Value = 6.85
$ws.Range("E777").Value = 7.28
$ws.Range("F777").Value = 64
$ws.Range("G777").Value = 438.4

$ws.Range("B784").Value = 46266
$ws.Range("D784").Value = 16.53
$ws.Range("E784").Value = 19.76
$ws.Range("F784").Value = 0
$ws.Range("G784").Value = 0

$ws.Range("B785").Value = 63807
$ws.Range("D785").Value = 16.53
$ws.Range("E785").Value = 17.58
$ws.Range("F785").Value = 96
$ws.Range("G785").Value = 1586.88

$ws.Range("B786").Value = 46268
$ws.Range("D786").Value = 26.38
$ws.Range("E786").Value = 31.51
$ws.Range("F786").Value = 0
$ws.Range("G786").Value = 0

$ws.Range("B787").Value = 63808
$ws.Range("D787").Value = 26.38
$ws.Range("E787").Value = 28.05
$ws.Range("F787").Value = 40
$ws.Range("G787").Value = 1055.2

$ws.Range("B807").Value = 64810
$ws.Range("D807").Value = 273.92
$ws.Range("E807").Value = 291.22
$ws.Range("F807").Value = 7
$ws.Range("G807").Value = 1917.44

$ws.Range("B808").Value = 53319
$ws.Range("D808").Value = 273.92
$ws.Range("E808").Value = 310.64
$ws.Range("F808").Value = 1
$ws.Range("G808").Value = 273.92

$ws.Range("B831").Value = 64832
$ws.Range("D831").Value = 32.83
$ws.Range("E831").Value = 34.9
$ws.Range("F831").Value = 100
$ws.Range("G831").Value = 3283

$ws.Range("B832").Value = 60024
$ws.Range("D832").Value = 32.83
$ws.Range("E832").Value = 37.22
$ws.Range("F832").Value = 0
$ws.Range("G832").Value = 0

$ws.Range("B833").Value = 60025
$ws.Range("D833").Value = 32.83
$ws.Range("E833").Value = 37.22
$ws.Range("F833").Value = 1
$ws.Range("G833").Value = 32.83

$ws.Range("B834").Value = 64833
$ws.Range("D834").Value = 32.83
$ws.Range("E834").Value = 34.9
$ws.Range("F834").Value = 99
$ws.Range("G834").Value = 3250.17

$ws.Range("B835").Value = 64831
$ws.Range("D835").Value = 32.83
$ws.Range("E835").Value = 34.9
$ws.Range("F835").Value = 152
$ws.Range("G835").Value = 4990.16

$ws.Range("B836").Value = 60023
$ws.Range("D836").Value = 32.83
$ws.Range("E836").Value = 37.22
$ws.Range("F836").Value = 0
$ws.Range("G836").Value = 0

$ws.Range("B839").Value = 60035
$ws.Range("D839").Value = 98.5
$ws.Range("E839").Value = 111.69
$ws.Range("F839").Value = 0
$ws.Range("G839").Value = 0

$ws.Range("B840").Value = 64839
$ws.Range("D840").Value = 98.5
$ws.Range("E840").Value = 104.71
$ws.Range("F840").Value = 94
$ws.Range("G840").Value = 9259

$ws.Range("B841").Value = 64838
$ws.Range("D841").Value = 98.5
$ws.Range("E841").Value = 104.71
$ws.Range("F841").Value = 78
$ws.Range("G841").Value = 7683

$ws.Range("B842").Value = 60034
$ws.Range("D842").Value = 98.5
$ws.Range("E842").Value = 111.69
$ws.Range("F842").Value = 0
$ws.Range("G842").Value = 0

$ws.Range("B853").Value = 60043
$ws.Range("D853").Value = 99.72
$ws.Range("E853").Value = 119.13
$ws.Range("F853").Value = 0
$ws.Range("G853").Value = 0

$ws.Range("B854").Value = 64843
$ws.Range("D854").Value = 99.72
$ws.Range("E854").Value = 106.01
$ws.Range("F854").Value = 2
$ws.Range("G854").Value = 199.44

$ws.Range("B870").Value = 54091
$ws.Range("D870").Value = 160.48
$ws.Range("E870").Value = 181.99
$ws.Range("F870").Value = 0
$ws.Range("G870").Value = 0

$ws.Range("B871").Value = 64173
$ws.Range("D871").Value = 160.48
$ws.Range("E871").Value = 170.61
$ws.Range("F871").Value = 2
$ws.Range("G871").Value = 320.96

$ws.Range("B902").Value = 41805
$ws.Range("D902").Value = 130.55
$ws.Range("E902").Value = 155.98
$ws.Range("F902").Value = 0
$ws.Range("G902").Value = 0

$ws.Range("B903").Value = 63736
$ws.Range("D903").Value = 130.55
$ws.Range("E903").Value = 138.78
$ws.Range("F903").Value = 37
$ws.Range("G903").Value = 4830.35

$ws.Range("B906").Value = 41800
$ws.Range("D906").Value = 130.55
$ws.Range("E906").Value = 155.98
$ws.Range("F906").Value = 0
$ws.Range("G906").Value = 0

$ws.Range("B907").Value = 63733
$ws.Range("D907").Value = 130.55
$ws.Range("E907").Value = 138.78
$ws.Range("F907").Value = 111
$ws.Range("G907").Value = 14491.05

$ws.Range("B939").Value = 61932
$ws.Range("D939").Value = 37783.32
$ws.Range("E939").Value = 48967.17
$ws.Range("F939").Value = 0
$ws.Range("G939").Value = 0

$ws.Range("B940").Value = 64614
$ws.Range("D940").Value = 37783.32
$ws.Range("E940").Value = 45141.63
$ws.Range("F940").Value = 2
$ws.Range("G940").Value = 75566.64

$ws.Range("B946").Value = 53475
$ws.Range("D946").Value = 577.29
$ws.Range("E946").Value = 654.65
$ws.Range("F946").Value = 0
$ws.Range("G946").Value = 0

$ws.Range("B947").Value = 64122
$ws.Range("D947").Value = 577.29
$ws.Range("E947").Value = 613.73
$ws.Range("F947").Value = 2
$ws.Range("G947").Value = 1154.58
